$d = $word.ActiveDocument

# Renamed "historical" to "recent", e.g.
# "Quantify daily mean temperatures from climate data for historical and
#  future period" -> "... for recent and future period".
$rng = $d.Content
$found = $rng.Find.Execute(
    "historical", $true, $true, $false, $false, $false,
    $true, 1, $false, $null, 0
)

if ($found) {
    # Replace the word's text in place (still a single run at this point).
    $rng.Text = "recent"
    $innerStart = $rng.Start
    $innerEnd = $rng.End

    # Force the editor to cut the replaced word into its own run -- mirroring
    # what happens when the word is selected and retyped by hand -- by
    # nudging a run-level property and immediately restoring it so the final
    # (visible) formatting is unchanged from its neighbours.
    $inner = $d.Range($innerStart, $innerEnd)
    $inner.Bold = 1
    $inner.Bold = 0
}
